$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.04"
$ws.Range("E2").Value = "'-3.46%"
$ws.Range("D3").Value = "'42.87"
$ws.Range("E3").Value = "'-6.56%"
$ws.Range("D4").Value = "'5.147"
$ws.Range("E4").Value = "'-9.23%"
$ws.Range("D5").Value = "'0.08163"
$ws.Range("E5").Value = "'-2.60%"
$ws.Range("D6").Value = "'4.327"
$ws.Range("E6").Value = "'-3.43%"
$ws.Range("E7").Value = "'-12.68%"
$ws.Range("D8").Value = "'0.9499"
$ws.Range("E8").Value = "'-4.33%"
$ws.Range("D9").Value = "'0.1113"
$ws.Range("E9").Value = "'-3.47%"
$ws.Range("D10").Value = "'0.1851"
$ws.Range("E10").Value = "'-4.05%"
$ws.Range("D11").Value = "'0.09403"
$ws.Range("E11").Value = "'-6.05%"
$ws.Range("D12").Value = "'0.04646"
$ws.Range("E12").Value = "'-0.63%"
$ws.Range("D13").Value = "'7.442"
$ws.Range("E13").Value = "'-28.10%"
$ws.Range("D14").Value = "'0.1058"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.001290"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("D16").Value = "'0.005942"
$ws.Range("E16").Value = "'-2.78%"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'-0.26%"
$ws.Range("E18").Value = "'-1.86%"
$ws.Range("E19").Value = "'-0.01%"
$ws.Range("E20").Value = "'-0.85%"
$ws.Range("D21").Value = "'0.2625"
$ws.Range("E21").Value = "'-1.05%"
$ws.Range("D22").Value = "'0.04184"
$ws.Range("E22").Value = "'-0.90%"
$ws.Range("D23").Value = "'0.001251"
$ws.Range("E23").Value = "'-4.44%"
$ws.Range("D24").Value = "'0.004298"
$ws.Range("E24").Value = "'-7.23%"
$ws.Range("D25").Value = "'0.0001112"
$ws.Range("E25").Value = "'-13.40%"
$ws.Range("D26").Value = "'0.0002981"
$ws.Range("E26").Value = "'-20.41%"
$ws.Range("D38").Value = "'0.02591"
$ws.Range("E38").Value = "'-7.36%"
$ws.Range("D39").Value = "'0.05537"
$ws.Range("E39").Value = "'-4.26%"
$ws.Range("D40").Value = "'0.007824"
$ws.Range("E40").Value = "'0.88%"
$ws.Range("D41").Value = "'0.1391"
$ws.Range("E41").Value = "'-3.04%"
$ws.Range("D42").Value = "'0.006606"
$ws.Range("E42").Value = "'-9.23%"
$ws.Range("D43").Value = "'0.002043"
$ws.Range("E43").Value = "'-3.29%"
$ws.Range("D44").Value = "'0.008447"
$ws.Range("E44").Value = "'-6.71%"
$ws.Range("D45").Value = "'0.3468"
$ws.Range("E45").Value = "'1.73%"
$ws.Range("D46").Value = "'0.00006989"
$ws.Range("E46").Value = "'-4.98%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.003481"
$ws.Range("E48").Value = "'-0.82%"
$ws.Range("D49").Value = "'0.003533"
$ws.Range("E49").Value = "'0.84%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.06%"
